# NAND and NOR 2 and 3 inputs
# Update the "Delta" transistor-width figures on the Inv1ma, Inv2mA and
# Inv3mA sheets (renamed project copied from Inversor -> NAND3), and move
# the active sheet / selections to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Inv1ma (2-input cell): I8 stops being a formula (=+I7/0.3) and becomes
# a plain entered value; I7 is re-entered too. Everything else on the
# sheet (N2, N3, J12, J13, J14) is a formula that depends on I8 and
# recalculates automatically.
# ---------------------------------------------------------------------
$wsInv1ma = $wb.Worksheets.Item("Inv1ma")
$wsInv1ma.Range("I7").Value = 37.8
$wsInv1ma.Range("I8").Value = 133

# ---------------------------------------------------------------------
# Inv2mA (3-input cell): only I8 changes, again dropping its formula in
# favour of a typed-in value.
# ---------------------------------------------------------------------
$wsInv2mA = $wb.Worksheets.Item("Inv2mA")
$wsInv2mA.Range("I8").Value = 87

# ---------------------------------------------------------------------
# Inv3mA: I8 becomes a typed value first (breaks the old I7<->I8 link),
# then I7 picks up a brand new formula driven off I8.
# ---------------------------------------------------------------------
$wsInv3mA = $wb.Worksheets.Item("Inv3mA")
$wsInv3mA.Range("I8").Value = 45
$wsInv3mA.Range("I7").Formula = "=+I8*0.3"

# ---------------------------------------------------------------------
# Selections / active sheet. Inv2mA and Inv3mA just move their cursor;
# Inv1ma becomes the new active tab, so activate + select it last.
# ---------------------------------------------------------------------
$wsInv2mA.Range("I9").Select()
$wsInv3mA.Range("E12").Select()

$wsInv1ma.Activate()
$wsInv1ma.Range("I9").Select()
